$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# --- Cell-level corrections in rows 2-24 ---
$ws.Cells.Item(2,6).Value = ""
$ws.Cells.Item(5,6).Value = 17.66
$ws.Cells.Item(6,3).Value = 15.1
$ws.Cells.Item(6,6).Value = 16.43
$ws.Cells.Item(8,3).Value = ""
$ws.Cells.Item(9,6).Value = ""
$ws.Cells.Item(10,6).Value = ""
$ws.Cells.Item(12,3).Value = 12.5
$ws.Cells.Item(14,3).Value = ""
$ws.Cells.Item(17,3).Value = 11.2
$ws.Cells.Item(18,3).Value = 11.5
$ws.Cells.Item(19,3).Value = ""
$ws.Cells.Item(20,3).Value = ""
$ws.Cells.Item(23,3).Value = 12.2
$ws.Cells.Item(24,6).Value = 16.78

# --- Full rewrite of rows 26-33 (data reshuffled / regenerated) ---
$ws.Cells.Item(26,1).Value = "SC 5"
$ws.Cells.Item(26,2).Value = -20.2
$ws.Cells.Item(26,3).Value = 10.8
$ws.Cells.Item(26,4).Value = -13.8
$ws.Cells.Item(26,5).Value = -5
$ws.Cells.Item(26,6).Value = 17.38
$ws.Cells.Item(27,1).Value = "SC 101"
$ws.Cells.Item(27,2).Value = -20.4
$ws.Cells.Item(27,3).Value = ""
$ws.Cells.Item(27,4).Value = -14.6
$ws.Cells.Item(27,5).Value = -10
$ws.Cells.Item(27,6).Value = 17
$ws.Cells.Item(28,1).Value = "SC 105"
$ws.Cells.Item(28,2).Value = ""
$ws.Cells.Item(28,3).Value = 11.1
$ws.Cells.Item(28,4).Value = -13.7
$ws.Cells.Item(28,5).Value = -5.9
$ws.Cells.Item(28,6).Value = ""
$ws.Cells.Item(29,1).Value = "SC 119"
$ws.Cells.Item(29,2).Value = ""
$ws.Cells.Item(29,3).Value = 11.2
$ws.Cells.Item(29,4).Value = -13
$ws.Cells.Item(29,5).Value = -6.8
$ws.Cells.Item(29,6).Value = 18.06
$ws.Cells.Item(30,1).Value = "SC 120"
$ws.Cells.Item(30,2).Value = -19.7
$ws.Cells.Item(30,3).Value = ""
$ws.Cells.Item(30,4).Value = -13.6
$ws.Cells.Item(30,5).Value = -5.7
$ws.Cells.Item(30,6).Value = 16.89
$ws.Cells.Item(31,1).Value = "SC 132"
$ws.Cells.Item(31,2).Value = -18.8
$ws.Cells.Item(31,3).Value = 15.3
$ws.Cells.Item(31,4).Value = -13.7
$ws.Cells.Item(31,5).Value = -8.1
$ws.Cells.Item(31,6).Value = 17.18
$ws.Cells.Item(32,1).Value = "SC 193"
$ws.Cells.Item(32,2).Value = ""
$ws.Cells.Item(32,3).Value = 10.5
$ws.Cells.Item(32,4).Value = -14.7
$ws.Cells.Item(32,5).Value = -6.4
$ws.Cells.Item(32,6).Value = 17.39
$ws.Cells.Item(33,1).Value = "SC 232"
$ws.Cells.Item(33,2).Value = -19.5
$ws.Cells.Item(33,3).Value = 10.4
$ws.Cells.Item(33,4).Value = -14.1
$ws.Cells.Item(33,5).Value = -10.7
$ws.Cells.Item(33,6).Value = 17.53

# --- Remove now-obsolete trailing rows 34 and 35 ---
$ws.Rows.Item(34).Delete()
$ws.Rows.Item(34).Delete()
